# Updated symbol list on Fri Dec 23 08:14:19 UTC 2022 with GitHub Actions
#
# Every data cell on this sheet (Price, Volume(1h), Hora, ...) is stored as literal
# text in the workbook, even when its content looks like a plain number (e.g. "7",
# "246.74"). Assigning a bare numeric-looking string via .Value would make Excel
# coerce it to a Number, which would change both the cell's type and its on-disk
# text representation (e.g. trailing zeros). To avoid that, numeric-looking
# replacement values are written with a leading quote character, which tells Excel
# to keep them as text ("quote prefix"); the cell's Style is then reset to "Normal"
# immediately afterwards so the quote-prefix text format isn't left behind and the
# cell keeps its original default formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = '''8'
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = '''21.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = '''8'
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = '''5.420'
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = '''8'
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = '''0.05793'
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = '''8'
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = '''3.387'
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = '''8'
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = '''6.327'
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = '''8'
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = '''0.8074'
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = '''8'
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = '''0.9580'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("G9").Value = '''8'
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.01122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("G10").Value = '''8'
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1428'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("G11").Value = '''8'
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07507'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("G12").Value = '''8'
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '''0.03239'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G13").Value = '''8'
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.03032'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("G14").Value = '''8'
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''4.152'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").Value = '''8'
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '''0.09400'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15BitMartTokenBMX'
$ws.Range("G16").Value = '''8'
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '''0.001607'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("G17").Value = '''8'
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '''0.04801'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("G18").Value = '''8'
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = '''0.006333'
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = '''8'
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = '''0.004096'
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Value = '''8'
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = '''0.0009970'
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Value = '''8'
$ws.Range("G21").Style = "Normal"
$ws.Range("G22").Value = '''8'
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = '''3.722'
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = '''8'
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = '''2.245'
$ws.Range("D24").Style = "Normal"
$ws.Range("G24").Value = '''8'
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = '''0.3236'
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = '''8'
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = '''0.1296'
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = '''8'
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = '''0.0003125'
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = '''8'
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = '''8'
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = '''8'
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = '''8'
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = '''8'
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = '''8'
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = '''8'
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = '''8'
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = '''8'
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = '''8'
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = '''8'
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = '''8'
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = '''8'
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = '''0.03887'
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = '''8'
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = '''0.006361'
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").Value = '''8'
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = '''0.1072'
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Value = '''8'
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = '''0.003001'
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").Value = '''8'
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = '''0.006697'
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = '''8'
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = '''0.00005593'
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = '''8'
$ws.Range("G45").Style = "Normal"
$ws.Range("G46").Value = '''8'
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = '''0.3901'
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = '''8'
$ws.Range("G47").Style = "Normal"
$ws.Range("G48").Value = '''8'
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = '''8'
$ws.Range("G49").Style = "Normal"
$ws.Range("G50").Value = '''8'
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = '''8'
$ws.Range("G51").Style = "Normal"
